$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously the last row (69) had the "date only" format to mark it as the
# latest entry. Now that we're appending row 70, row 69 reverts to the
# regular "date + time" format and row 70 becomes the new "latest" row.
$ws.Range("A69").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A70").Value = 45810
$ws.Range("A70").NumberFormat = "YYYY-MM-DD"
$ws.Range("B70").Value = 299
$ws.Range("C70").Value = 294
$ws.Range("D70").Value = 298
